$wb = $excel.ActiveWorkbook

# Insert a new worksheet right after "Sheet1", named "new_setting"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "new_setting"

# --- Header row 3 (same as Sheet1) ---
$newSheet.Range("D3").Value = "K"
$newSheet.Range("E3").Value = "n_i"
$newSheet.Range("F3").Value = "lamba"
$newSheet.Range("G3").Value = "response"
$newSheet.Range("H3").Value = "b1"
$newSheet.Range("I3").Value = "b2"
$newSheet.Range("J3").Value = "price"
$newSheet.Range("K3").Value = "T_D"
$newSheet.Range("L3").Value = "Vd(Mb)"
$newSheet.Range("M3").Value = "cost"
$newSheet.Range("N3").Value = "t_l"
$newSheet.Range("O3").Value = "packet size"
$newSheet.Range("P3").Value = "ping time"

# --- Table 1: rows 4-8 (identical values/formulas to Sheet1) ---
$aps = @("AP1","AP2","AP3","AP4","AP5")
$d = @(3,3,3,3,3)
$e = @(2,4,10,5,15)
$h = @(5,5,5,5,5)
$i = @(5,10,30,10,15)
$j = @(0.5,0.75,1.25,0.5,1)
$k = @(100,100,100,100,100)
$l = @(250,250,250,250,250)
$o = @(2,2,2,2,2)
for ($r = 4; $r -le 8; $r++) {
    $idx = $r - 4
    $newSheet.Range("C$r").Value = $aps[$idx]
    $newSheet.Range("D$r").Value = $d[$idx]
    $newSheet.Range("E$r").Value = $e[$idx]
    $newSheet.Range("F$r").Formula = "=1/P$r"
    $newSheet.Range("G$r").Formula = "=E$r/F$r"
    $newSheet.Range("H$r").Value = $h[$idx]
    $newSheet.Range("I$r").Value = $i[$idx]
    $newSheet.Range("J$r").Value = $j[$idx]
    $newSheet.Range("K$r").Value = $k[$idx]
    $newSheet.Range("L$r").Value = $l[$idx]
    $newSheet.Range("M$r").Formula = "=L$r*J$r"
    $newSheet.Range("N$r").Formula = "=(L$r*8)/(H$r+I$r)"
    $newSheet.Range("O$r").Value = $o[$idx]
    $newSheet.Range("P$r").Formula = "=(O$r*8)/(H$r+I$r)"
}

# --- Row 11 header (Ranking table) ---
$newSheet.Range("B11").Value = "Ranking"
$newSheet.Range("D11").Value = "K"
$newSheet.Range("E11").Value = "n_i"
$newSheet.Range("F11").Value = "lamba"
$newSheet.Range("G11").Value = "response"
$newSheet.Range("H11").Value = "b1"
$newSheet.Range("I11").Value = "b2"
$newSheet.Range("J11").Value = "price"
$newSheet.Range("K11").Value = "T_D"
$newSheet.Range("L11").Value = "Vd(Mb)"
$newSheet.Range("M11").Value = "cost"
$newSheet.Range("N11").Value = "t_l"
$newSheet.Range("O11").Value = "packet size"
$newSheet.Range("P11").Value = "ping time"

# --- Table 2: rows 12-17 (Ranking results; I column now uniformly 60) ---
$rank = @(1,2,4,3,5,$null)
$aps2 = @("AP1","AP2","AP3","AP4","AP5","AP")
$d2 = @(3,3,3,3,3,3)
$e2 = @(2,4,10,5,15,30)
$h2 = @(20,20,20,20,20,20)
$i2 = @(60,60,60,60,60,60)
$j2 = @(0.02,0.05,0.01,0.02,0.01,0)
$k2 = @(500,500,500,500,500,500)
$l2 = @(3072,3072,3072,3072,3072,3072)
$o2 = @(2,2,2,2,2,2)
for ($r = 12; $r -le 17; $r++) {
    $idx = $r - 12
    if ($rank[$idx] -ne $null) {
        $newSheet.Range("B$r").Value = $rank[$idx]
    }
    $newSheet.Range("C$r").Value = $aps2[$idx]
    $newSheet.Range("D$r").Value = $d2[$idx]
    $newSheet.Range("E$r").Value = $e2[$idx]
    $newSheet.Range("F$r").Formula = "=1/P$r"
    $newSheet.Range("G$r").Formula = "=E$r/F$r"
    $newSheet.Range("H$r").Value = $h2[$idx]
    $newSheet.Range("I$r").Value = $i2[$idx]
    $newSheet.Range("J$r").Value = $j2[$idx]
    $newSheet.Range("K$r").Value = $k2[$idx]
    $newSheet.Range("L$r").Value = $l2[$idx]
    $newSheet.Range("M$r").Formula = "=L$r*J$r"
    $newSheet.Range("N$r").Formula = "=(L$r*8)/(H$r+I$r)"
    $newSheet.Range("O$r").Value = $o2[$idx]
    $newSheet.Range("P$r").Formula = "=(O$r*8)/(H$r+I$r)"
}

$newSheet.Range("M18").Formula = "=AVERAGE(M12:M16)"

# --- Row 19 header (b2 sim table) ---
$newSheet.Range("D19").Value = "b2"
$newSheet.Range("E19").Value = "r1"
$newSheet.Range("F19").Value = "r2"
$newSheet.Range("G19").Value = "r3"
$newSheet.Range("H19").Value = "t1"
$newSheet.Range("I19").Value = "t2"
$newSheet.Range("J19").Value = "t3"
$newSheet.Range("K19").Value = "E[b2]"

# --- Table 3: rows 20-25 (new data values, different from Sheet1) ---
$aps3 = @("AP1","AP2","AP3","AP4","AP5","AP6")
$d3 = @(60,60,60,60,60,60)
$e3 = @(40,30,5,5,5,1)
$f3 = @(60,50,80,50,80,65)
$g3 = @(80,100,90,120,90,100)
$h3 = @(30,60,150,100,5,300)
$i3 = @(30,60,150,100,5,300)
$j3 = @(30,60,150,100,5,300)
for ($r = 20; $r -le 25; $r++) {
    $idx = $r - 20
    $newSheet.Range("C$r").Value = $aps3[$idx]
    $newSheet.Range("D$r").Value = $d3[$idx]
    $newSheet.Range("E$r").Value = $e3[$idx]
    $newSheet.Range("F$r").Value = $f3[$idx]
    $newSheet.Range("G$r").Value = $g3[$idx]
    $newSheet.Range("H$r").Value = $h3[$idx]
    $newSheet.Range("I$r").Value = $i3[$idx]
    $newSheet.Range("J$r").Value = $j3[$idx]
    $newSheet.Range("K$r").Formula = "=((E$r*H$r)+(F$r*I$r)+(G$r*J$r))/(H$r+I$r+J$r)"
}

$newSheet.Range("E20").Select() | Out-Null

$sheet1.Activate() | Out-Null
$sheet1.Range("N12:N16").Select() | Out-Null
